# Weekly update: add a new price record for "Pepino dulce" at the top of
# the data block (new row 19), pushing all existing records (old rows
# 19-55) down by one row to rows 20-56.
#
# The new row 19 reuses the same boilerplate columns (Mercado, Region,
# Codreg, Categoria ID/Nombre, Variedad, Unidad de comercializacion,
# Origen, Kg o Unidades, Clasificacion) as the record that is now in row
# 20 (these are constant across every record in this sheet), and carries
# its own Fecha / Calidad / Volumen values while keeping the same
# Precio minimo/maximo/promedio ponderado and Precio $/Kg as before.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at row 19; rows 19-55 shift down to 20-56.
$ws.Rows.Item(19).Insert()

# Duplicate the (now shifted) row 20 into the new row 19 for all 18
# columns (A:R) so the boilerplate fields line up.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(19, $col).Value = $ws.Cells.Item(20, $col).Value2
}

# Overwrite the fields that differ for this new weekly record.
$ws.Cells.Item(19, 4).Value = 44662    # D19 - Fecha
$ws.Cells.Item(19, 9).Value = "Primera" # I19 - Calidad
$ws.Cells.Item(19, 10).Value = 50       # J19 - Volumen
